$d = $word.ActiveDocument

# Land a fresh paragraph break at the very end of the document (after the
# "...Otherwise the highest versioned SDK..." paragraph), then replace that
# new (still-empty) paragraph's range -- mark included -- with the whole
# new block of content in one shot via InsertXML. Because the target range
# we replace includes the paragraph mark, Word preserves a trailing empty
# paragraph (with the same pPr) after the inserted content, which gives us
# the blank line that closes out the document.
$tail = $d.Range($d.Content.End, $d.Content.End)
$tail.InsertParagraphAfter()

$newPara = $d.Paragraphs($d.Paragraphs.Count)
$target = $newPara.Range

$xml = '<?xml version="1.0" encoding="UTF-8" standalone="yes"?>' +
'<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage">' +
'<pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml">' +
'<pkg:xmlData>' +
'<w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main">' +
'<w:body>' +
  '<w:p><w:pPr><w:jc w:val="both"/></w:pPr></w:p>' +
  '<w:p><w:pPr><w:jc w:val="both"/><w:rPr><w:b/><w:bCs/></w:rPr></w:pPr><w:r><w:rPr><w:b/><w:bCs/></w:rPr><w:t>Echo</w:t></w:r></w:p>' +
  '<w:p><w:pPr><w:jc w:val="both"/></w:pPr>' +
    '<w:r><w:t>Echo "</w:t></w:r>' +
    '<w:proofErr w:type="spellStart"/>' +
    '<w:r><w:t>JupyterLab</w:t></w:r>' +
    '<w:proofErr w:type="spellEnd"/>' +
    '<w:r><w:t xml:space="preserve"> is up and running" into a text file named jupyter_test.txt:</w:t></w:r>' +
  '</w:p>' +
  '<w:p><w:pPr><w:jc w:val="both"/></w:pPr>' +
    '<w:r><w:t>echo ' + [char]39 + '</w:t></w:r>' +
    '<w:proofErr w:type="spellStart"/>' +
    '<w:r><w:t>Jupyter</w:t></w:r>' +
    '<w:proofErr w:type="spellEnd"/>' +
    '<w:r><w:t xml:space="preserve"> is up and running' + [char]39 + ' &gt; jupyter_test.txt</w:t></w:r>' +
  '</w:p>' +
'</w:body>' +
'</w:document>' +
'</pkg:xmlData>' +
'</pkg:part>' +
'</pkg:package>'

$target.InsertXML($xml)

# Drop the now-obsolete <w:semiHidden/> flag on the DefaultParagraphFont
# character style.
$styles = $d.Styles
$dpf = $styles.Item("Default Paragraph Font")
$dpf.SemiHidden = $false
